$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.672615647315979
$ws.Range("B1").Value = 2.492033243179321
$ws.Range("C1").Value = 2.864604234695435
$ws.Range("D1").Value = 3.304876089096069
$ws.Range("E1").Value = 1.046448111534119
